$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column N (shifts Late/heading/Outstanding
# columns from N,O,P to O,P,Q), matching column M's width for the new column.
$mWidth = $ws.Columns("M").ColumnWidth
$ws.Columns("N").Insert()
$ws.Columns("N").ColumnWidth = $mWidth

# Make "Repayment schedule" the active sheet/tab with the given selection.
[void]$ws.Activate()
[void]$ws.Range("L13").Select()
